$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.567823529243469
$ws.Range("B1").Value = 1.733191132545471
$ws.Range("C1").Value = 1.758923411369324
$ws.Range("D1").Value = 2.33620810508728
$ws.Range("E1").Value = 3.979049921035767
